function Set-ParagraphXml($d, $index, $innerXml) {
    $p = $d.Paragraphs($index)
    $full = $p.Range
    $content = $d.Range($full.Start, $full.End - 1)
    $pkg = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p>" + $innerXml + "</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    $null = $content.InsertXML($pkg)
}

$d = $word.ActiveDocument

$xml1 = '<w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00B42414"><w:rPr><w:rFonts w:ascii="Segoe Print" w:hAnsi="Segoe Print"/><w:sz w:val="96"/><w:szCs w:val="96"/></w:rPr><w:t>Шестёрочка</w:t></w:r><w:proofErr w:type="spellEnd"/>'
Set-ParagraphXml $d 1 $xml1

$xml3 = '<w:r w:rsidRPr="00B42414"><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Цена</w:t></w:r><w:r w:rsidR="0095142F" w:rsidRPr="0095142F"><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">             </w:t></w:r><w:r w:rsidRPr="00B42414"><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Скидка</w:t></w:r><w:r w:rsidR="0095142F" w:rsidRPr="0095142F"><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">            </w:t></w:r><w:r w:rsidRPr="00B42414"><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Кол-во</w:t></w:r><w:r w:rsidR="0095142F" w:rsidRPr="0095142F"><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">           </w:t></w:r><w:r w:rsidR="007D7754" w:rsidRPr="00B42414"><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Итого</w:t></w:r>'
Set-ParagraphXml $d 3 $xml3

$xml5 = '<w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00B42414"><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">{{ </w:t></w:r><w:r w:rsidR="005010D6"><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t>item</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00B42414"><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t>}}</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
Set-ParagraphXml $d 5 $xml5

$xml6 = '<w:r w:rsidRPr="00B42414"><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00B42414"><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t>endfor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00B42414"><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> %}</w:t></w:r>'
Set-ParagraphXml $d 6 $xml6

$xml8 = '<w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>ИТОГ</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t>{{ total</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> }}</w:t></w:r>'
Set-ParagraphXml $d 8 $xml8

$xml9 = '<w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>СКИДКА</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t>{{ discount</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> }}</w:t></w:r>'
Set-ParagraphXml $d 9 $xml9

$xml11 = '<w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Кассир</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t>:</w:t></w:r><w:r w:rsidR="00E20DDD"><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">   </w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="00E20DDD"><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">       </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t>{{ cashier }}</w:t></w:r>'
Set-ParagraphXml $d 11 $xml11

$xml12 = '<w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t>{{ date</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> }}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiLight" w:hAnsi="Bahnschrift SemiLight"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t>{{ time }}</w:t></w:r>'
Set-ParagraphXml $d 12 $xml12

Write-Output "All paragraphs updated."